$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 160
$ws.Range("I18").Value = 160
$ws.Range("K18").Value = 160
$ws.Range("M18").Value = 124
$ws.Range("H62").Value = 3518.9565
$ws.Range("I62").Value = 2926.6667
$ws.Range("J62").Value = 4629.5
$ws.Range("K62").Value = 2926.6667
$ws.Range("L62").Value = 4629.5
$ws.Range("M62").Value = -2302.6667
$ws.Range("N62").Value = -5877.5
$ws.Range("H64").Value = 3468.9795
$ws.Range("I64").Value = 3191.795
$ws.Range("J64").Value = 4550
$ws.Range("K64").Value = 3191.795
$ws.Range("L64").Value = 4550
$ws.Range("M64").Value = -2943.795
$ws.Range("N64").Value = -5046
$ws.Range("H65").Value = 3518.9565
$ws.Range("I65").Value = 2926.6667
$ws.Range("J65").Value = 4629.5
$ws.Range("K65").Value = 14633.3335
$ws.Range("L65").Value = 23147.5
$ws.Range("M65").Value = -11513.3335
$ws.Range("N65").Value = -29387.5
$ws.Range("H67").Value = 3468.9795
$ws.Range("I67").Value = 3191.795
$ws.Range("J67").Value = 4550
$ws.Range("K67").Value = 3191.795
$ws.Range("L67").Value = 4550
$ws.Range("M67").Value = -2333.795
$ws.Range("N67").Value = -6266
$ws.Range("H76").Value = 3653.4565
$ws.Range("I76").Value = 3001.4707
$ws.Range("J76").Value = 5500.75
$ws.Range("K76").Value = 3001.4707
$ws.Range("L76").Value = 5500.75
$ws.Range("M76").Value = -2686.4707
$ws.Range("N76").Value = -6130.75
$ws.Range("H79").Value = 3653.4565
$ws.Range("I79").Value = 3001.4707
$ws.Range("J79").Value = 5500.75
$ws.Range("K79").Value = 3001.4707
$ws.Range("L79").Value = 5500.75
$ws.Range("M79").Value = -1909.4707
$ws.Range("N79").Value = -7684.75
$ws.Range("H86").Value = 56146.637
$ws.Range("J86").Value = 135268
$ws.Range("L86").Value = 135268
$ws.Range("N86").Value = -137514
$ws.Range("H87").Value = 12339.23
$ws.Range("J87").Value = 12339.23
$ws.Range("L87").Value = 12339.23
$ws.Range("N87").Value = -14835.23
$ws.Range("H89").Value = 56146.637
$ws.Range("J89").Value = 135268
$ws.Range("L89").Value = 676340
$ws.Range("N89").Value = -687572
$ws.Range("H90").Value = 12339.23
$ws.Range("J90").Value = 12339.23
$ws.Range("L90").Value = 37017.69
$ws.Range("N90").Value = -49497.69
$ws.Range("H98").Value = 1736.3125
$ws.Range("I98").Value = 2549.75
$ws.Range("J98").Value = 1465.1666
$ws.Range("K98").Value = 2549.75
$ws.Range("L98").Value = 1465.1666
$ws.Range("M98").Value = -1051.75
$ws.Range("N98").Value = -4461.1666
$ws.Range("H122").Value = 1736.3125
$ws.Range("I122").Value = 2549.75
$ws.Range("J122").Value = 1465.1666
$ws.Range("K122").Value = 7649.25
$ws.Range("L122").Value = 4395.4998
$ws.Range("M122").Value = -5199.25
$ws.Range("N122").Value = -9295.4998
$ws.Range("H124").Value = 12000
$ws.Range("J124").Value = 12000
$ws.Range("L124").Value = 12000
$ws.Range("N124").Value = -21820
$ws.Range("H126").Value = 11995.883
$ws.Range("J126").Value = 11995.883
$ws.Range("L126").Value = 11995.883
$ws.Range("N126").Value = -21875.883
$ws.Range("H127").Value = 2062.6667
$ws.Range("J127").Value = 2241.257
$ws.Range("L127").Value = 6723.771000000001
$ws.Range("N127").Value = -16643.771
$ws.Range("H128").Value = 11994.546
$ws.Range("J128").Value = 11994.546
$ws.Range("L128").Value = 11994.546
$ws.Range("N128").Value = -21954.546
$ws.Range("H130").Value = 13995.77
$ws.Range("J130").Value = 13995.77
$ws.Range("L130").Value = 13995.77
$ws.Range("N130").Value = -24035.77
$ws.Range("H133").Value = 44000
$ws.Range("J133").Value = 44000
$ws.Range("L133").Value = 44000
$ws.Range("N133").Value = -54120
$ws.Range("H138").Value = 3341.5408
$ws.Range("I138").Value = 1588.0952
$ws.Range("K138").Value = 4764.2856
$ws.Range("M138").Value = 375.7143999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6878.6865
$ws.Range("I32").Value = 5542.18
$ws.Range("J32").Value = 20466.5
$ws.Range("K32").Value = 5542.18
$ws.Range("L32").Value = 20466.5
$ws.Range("M32").Value = -5255.18
$ws.Range("N32").Value = -21040.5
$ws.Range("H61").Value = 2080.8333
$ws.Range("I61").Value = 1949.5238
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1949.5238
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1737.5238
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 4460.923
$ws.Range("I74").Value = 5785.3335
$ws.Range("K74").Value = 5785.3335
$ws.Range("M74").Value = -4911.3335
$ws.Range("H77").Value = 4460.923
$ws.Range("I77").Value = 5785.3335
$ws.Range("K77").Value = 28926.6675
$ws.Range("M77").Value = -24558.6675
$ws.Range("H97").Value = 720.7143
$ws.Range("I97").Value = 632.5
$ws.Range("J97").Value = 838.3333
$ws.Range("K97").Value = 632.5
$ws.Range("L97").Value = 838.3333
$ws.Range("M97").Value = -136.5
$ws.Range("N97").Value = -1830.3333
$ws.Range("H132").Value = 1303.1017
$ws.Range("I132").Value = 788.4524
$ws.Range("J132").Value = 2574.5881
$ws.Range("K132").Value = 2365.3572
$ws.Range("L132").Value = 7723.7643
$ws.Range("M132").Value = 164.6428000000001
$ws.Range("N132").Value = -12783.7643
$ws.Range("H136").Value = 2080.8333
$ws.Range("I136").Value = 1949.5238
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5848.5714
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3298.5714
$ws.Range("N136").Value = -14100
$ws.Range("H139").Value = 38716.152
$ws.Range("J139").Value = 38716.152
$ws.Range("L139").Value = 38716.152
$ws.Range("N139").Value = -48996.152

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8833.267
$ws.Range("I20").Value = 10183.25
$ws.Range("J20").Value = 3433.3333
$ws.Range("K20").Value = 10183.25
$ws.Range("L20").Value = 3433.3333
$ws.Range("M20").Value = -9936.25
$ws.Range("N20").Value = -3927.3333
$ws.Range("H94").Value = 628
$ws.Range("I94").Value = 723.6
$ws.Range("K94").Value = 723.6
$ws.Range("M94").Value = -272.6
$ws.Range("H134").Value = 1406.8368
$ws.Range("I134").Value = 1203.05
$ws.Range("J134").Value = 2312.5557
$ws.Range("K134").Value = 3609.15
$ws.Range("L134").Value = 6937.6671
$ws.Range("M134").Value = -1074.15
$ws.Range("N134").Value = -12007.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2745.5
$ws.Range("I22").Value = 2745.5
$ws.Range("K22").Value = 2745.5
$ws.Range("M22").Value = -2395.5
$ws.Range("H86").Value = 5348.5625
$ws.Range("I86").Value = 5744.4443
$ws.Range("J86").Value = 4839.5713
$ws.Range("K86").Value = 5744.4443
$ws.Range("L86").Value = 4839.5713
$ws.Range("M86").Value = -4621.4443
$ws.Range("N86").Value = -7085.5713
$ws.Range("H89").Value = 5348.5625
$ws.Range("I89").Value = 5744.4443
$ws.Range("J89").Value = 4839.5713
$ws.Range("K89").Value = 28722.2215
$ws.Range("L89").Value = 24197.8565
$ws.Range("M89").Value = -23106.2215
$ws.Range("N89").Value = -35429.85649999999
$ws.Range("H132").Value = 2765.4375
$ws.Range("I132").Value = 2261.6843
$ws.Range("J132").Value = 3501.6924
$ws.Range("K132").Value = 6785.0529
$ws.Range("L132").Value = 10505.0772
$ws.Range("M132").Value = -4255.0529
$ws.Range("N132").Value = -15565.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 18753.666
$ws.Range("J39").Value = 18753.666
$ws.Range("L39").Value = 18753.666
$ws.Range("N39").Value = -19817.666
$ws.Range("H70").Value = 5493.375
$ws.Range("I70").Value = 4660.0435
$ws.Range("K70").Value = 4660.0435
$ws.Range("M70").Value = -4390.0435
$ws.Range("H73").Value = 5493.375
$ws.Range("I73").Value = 4660.0435
$ws.Range("K73").Value = 4660.0435
$ws.Range("M73").Value = -3724.0435
$ws.Range("H80").Value = 2165
$ws.Range("J80").Value = 2495
$ws.Range("L80").Value = 2495
$ws.Range("N80").Value = -4491
$ws.Range("H83").Value = 2165
$ws.Range("J83").Value = 2495
$ws.Range("L83").Value = 12475
$ws.Range("N83").Value = -22459

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 17599832
$ws.Range("I136").Value = 80382.38
$ws.Range("J136").Value = 55558640
$ws.Range("K136").Value = 241147.14
$ws.Range("L136").Value = 166675920
$ws.Range("M136").Value = -238597.14
$ws.Range("N136").Value = -166681020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1640
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 1850
$ws.Range("K96").Value = 800
$ws.Range("L96").Value = 1850
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -4596
